$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.174.84"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.329.96"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "253.85"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("D7").Value = "75.93"
$ws.Range("E7").Value = "  +5.96%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "0.658"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("D10").Value = "40.49"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").Value = "0.0990"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "7.58"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "2.673.10"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "15.49"
$ws.Range("E15").Value = "  +4.00%  "
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "2.327.28"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "43.147.73"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").Value = "6.35"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "73.19"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "239.13"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("E23").Value = "  +5.66%  "
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "11.69"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "21.35"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "167.60"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "6.35"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "0.0851"
$ws.Range("E32").Value = "  +10.11%  "
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Value = "30.61"
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("D35").Value = "0.128"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "4.62"
$ws.Range("E36").Value = "  +10.93%  "
$ws.Range("D37").Value = "4.87"
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").Value = "13.93"
$ws.Range("E39").Value = "  +13.52%  "
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "5.95"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").Value = "0.220"
$ws.Range("E42").Value = "  +8.59%  "
$ws.Range("D43").Value = "9.26"
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("D44").Value = "62.89"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("D45").Value = "4.92"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("D46").Value = "106.03"
$ws.Range("E46").Value = "  +11.47%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D50").Value = "1.20"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "4.37"
$ws.Range("E51").Value = "  -0.41%  "
